# Resubmission version with data
# Refresh the regression-table values with the re-run's (numerically
# near-identical) output, and clear the stray duplicate "N (villages)"-style
# row (row 30), which had repeated the row-29 values by mistake.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated coefficient / SE values (tiny re-estimation deltas) ---
$ws.Range("C2").Value  = 22.542368594070389
$ws.Range("D2").Value  = 22.845221309041861
$ws.Range("C3").Value  = 4.9914988197683501
$ws.Range("D4").Value  = 8.1513283337240523
$ws.Range("B5").Value  = 2.9558869297900752
$ws.Range("C5").Value  = 3.0211278568390929
$ws.Range("B8").Value  = 0.97533653023810207
$ws.Range("B10").Value = -1.087693351042927
$ws.Range("C10").Value = -1.1686110092835129
$ws.Range("D10").Value = -0.97854815528164496
$ws.Range("B11").Value = 2.9786569508931442
$ws.Range("C11").Value = 3.1236596236980931
$ws.Range("D11").Value = 3.3155023501835612
$ws.Range("B12").Value = 1.3838382209480931
$ws.Range("C12").Value = 1.3205404370178599
$ws.Range("D12").Value = 1.3588563053771141
$ws.Range("B13").Value = 2.9852350131827898
$ws.Range("C13").Value = 3.132881995205941
$ws.Range("D13").Value = 3.3162153979151392
$ws.Range("B14").Value = 2.3836781147719579
$ws.Range("C14").Value = 2.5189349215756822
$ws.Range("D14").Value = 2.806064705483164
$ws.Range("B15").Value = 3.0121825195861738
$ws.Range("C15").Value = 3.1622207095883832
$ws.Range("D15").Value = 3.344866158287199
$ws.Range("C17").Value = 0.488271746730923
$ws.Range("D17").Value = 0.50541704621996697
$ws.Range("B18").Value = 0.94477636554091027
$ws.Range("D18").Value = 0.65171488956250656
$ws.Range("C19").Value = 0.6414990894536835
$ws.Range("D19").Value = 0.63462036668167543
$ws.Range("B20").Value = -0.99314719384658878
$ws.Range("C20").Value = -0.84521752848743925
$ws.Range("D20").Value = -1.2148272520982519
$ws.Range("D21").Value = 1.164011047580247
$ws.Range("C23").Value = 0.90974262819674945
$ws.Range("C25").Value = 0.018081189215723801
$ws.Range("B26").Value = 3.3912368556378159
$ws.Range("C26").Value = -0.1322745539828496
$ws.Range("D26").Value = 3.365866315888054
$ws.Range("B27").Value = 3.2625453487766229
$ws.Range("C27").Value = 3.654456034219733
$ws.Range("D27").Value = 3.5659531407363541

# --- Clear the duplicated row 30 (kept the label, dropped the stray values) ---
$ws.Range("B30:D30").ClearContents()
